$d = $word.ActiveDocument

# Before the edit, the document ends with:
#   <w:p/>
#   <w:p/>
#   <w:p> (underlined) "Todos los primitivos se pasan por valor, y todos los objetos por refencia."
#
# After the edit it must read:
#   <w:p/>
#   <w:p/>
#   <w:p> (plain)      "Todos los primitivos se pasan por valor, y todos los objetos por refencia."
#   <w:p/>
#   <w:p> (underlined) "En los proyectos en la carpeta assests suele ser para recursos estáticos que no van a ser eliminados."
#
# i.e. the original sentence is moved up into a brand new, unformatted paragraph
# (followed by a fresh blank paragraph), and the still-underlined paragraph that
# used to hold it gets a brand new sentence instead.

$oldSentence = "Todos los primitivos se pasan por valor, y todos los objetos por refencia."
$newSentence = "En los proyectos en la carpeta assests suele ser para recursos estáticos que no van a ser eliminados."

# Locate the (currently last) underlined paragraph that holds $oldSentence.
$underlinedIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "$oldSentence*") {
        $underlinedIndex = $i
    }
}

# The paragraph right before it is one of the existing blank paragraphs. Inserting
# a new paragraph after a plain/blank paragraph produces a clean, unformatted new
# paragraph (it does not inherit the underline formatting of the paragraph that
# follows it).
$precedingPara = $d.Paragraphs.Item($underlinedIndex - 1)
$precedingPara.Range.InsertParagraphAfter()

# The sentence moves into this freshly inserted (plain) paragraph.
$movedPara = $d.Paragraphs.Item($underlinedIndex - 1 + 1)
$movedPara.Range.Text = $oldSentence

# A new blank paragraph separates the moved sentence from the underlined paragraph.
$movedPara = $d.Paragraphs.Item($underlinedIndex - 1 + 1)
$movedPara.Range.InsertParagraphAfter()

# Finally, swap the text inside the still-underlined paragraph (it kept its
# formatting and simply shifted two paragraphs further down).
$underlinedPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$underlinedPara.Range.Find.Execute($oldSentence, $true, $false, $false, $false, $false, `
    $true, 1, $false, $newSentence, 2)
